# Apply the "Drools integrated and running code" edit to the rule.xlsx
# decision table:
#   - C4 ("Sequential" row) switches from a boolean TRUE to the literal
#     text string "true" (quotes included), rendered in a slightly
#     smaller, left-aligned font.
#   - C2 ("RuleSet" row) changes from "rules" to "com.livelabdrools".
#   - The active selection moves from C3 to C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- C4: Sequential value becomes the literal string "true" ---------------
$c4 = $ws.Range("C4")
$c4.ClearContents()
$c4.Value = '"true"'
$c4.Font.Size = 10
$c4.HorizontalAlignment = -4131   # xlHAlignLeft

# --- C2: RuleSet value becomes com.livelabdrools ---------------------------
$ws.Range("C2").Value = "com.livelabdrools"

# --- update the selected cell to C2 ----------------------------------------
[void]$ws.Range("C2").Select()
